$wb = $excel.ActiveWorkbook

# --- POBasedInvoice sheet: refresh test invoice data (rows 2-9) ---
$po = $wb.Worksheets.Item("POBasedInvoice")

# Columns A (Invoice Number), B (Base Amount) and C (IGST) are stored as
# text even though their values look numeric, so force text formatting
# before writing, then drop back to the Normal style (matches the
# source file, which keeps these cells on the default style).
$po.Range("A2:C9").NumberFormat = "@"

$po.Range("A2").Value = "TESTINV3712"
$po.Range("B2").Value = "1"
$po.Range("C2").Value = "0.18"

$po.Range("A3").Value = "TESTINV1887"
$po.Range("B3").Value = "7"
$po.Range("C3").Value = "1.26"

$po.Range("A4").Value = "TESTINV5595"
$po.Range("B4").Value = "4"
$po.Range("C4").Value = "0.72"

$po.Range("A5").Value = "TESTINV2223"
$po.Range("B5").Value = "2"
$po.Range("C5").Value = "0.36"

$po.Range("A6").Value = "TESTINV2561"
$po.Range("B6").Value = "4"
$po.Range("C6").Value = "0.72"

$po.Range("A7").Value = "TESTINV4211"
$po.Range("B7").Value = "6"
$po.Range("C7").Value = "1.08"

$po.Range("A8").Value = "TESTINV3819"
$po.Range("B8").Value = "8"
$po.Range("C8").Value = "1.44"

$po.Range("A9").Value = "TESTINV7271"
$po.Range("B9").Value = "6"
$po.Range("C9").Value = "1.08"

$po.Range("A2:C9").Style = "Normal"

# Quantity column drops from 10 to 1 for every invoice row (numeric cells).
$po.Range("J2").Value = 1
$po.Range("J3").Value = 1
$po.Range("J4").Value = 1
$po.Range("J5").Value = 1
$po.Range("J6").Value = 1
$po.Range("J7").Value = 1
$po.Range("J8").Value = 1
$po.Range("J9").Value = 1

# --- BADashboardPage sheet: update "To state" value ---
$dash = $wb.Worksheets.Item("BADashboardPage")
$dash.Range("B2").Value = "MANIPUR"

# --- Make POBasedInvoice the active sheet/tab with A2 selected ---
$po.Activate()
$po.Range("A2").Select()
